$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 236).
$ws.Range("C2:C236").Value = 45180
